$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-10-04 01:10:48"

for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
